$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 195.435389
$ws.Range("H2").Value = 586.306167
$ws.Range("I2").Value = 0.3095741734129938
$ws.Range("J2").Value = 0.3095741734129938
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.002414333333333
$ws.Range("N2").Value = 6.007243
$ws.Range("O2").Value = 0.8176262304899502
$ws.Range("P2").Value = 0.8176262304899502
$ws.Range("Q2").Value = 391.3426241741756
$ws.Range("R2").Value = 3522.083617567581
$ws.Range("S2").Value = 0.2531159644647083
$ws.Range("T2").Value = 0.2531159644647082

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 195.435389
$ws.Range("H3").Value = 586.306167
$ws.Range("I3").Value = 0.3095741734129938
$ws.Range("J3").Value = 0.3095741734129938
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.446644
$ws.Range("N3").Value = 1.339932
$ws.Range("O3").Value = 0.1823737695100498
$ws.Range("P3").Value = 0.1823737695100498
$ws.Range("Q3").Value = 87.290043884516
$ws.Range("R3").Value = 785.610394960644
$ws.Range("S3").Value = 0.05645820894828552
$ws.Range("T3").Value = 0.05645820894828551

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 79.82725266666667
$ws.Range("H4").Value = 239.481758
$ws.Range("I4").Value = 0.1264482133280045
$ws.Range("J4").Value = 0.1264482133280045
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.002414333333333
$ws.Range("N4").Value = 6.007243
$ws.Range("O4").Value = 0.8176262304899502
$ws.Range("P4").Value = 0.8176262304899502
$ws.Range("Q4").Value = 159.8472349303549
$ws.Range("R4").Value = 1438.625114373194
$ws.Range("S4").Value = 0.1033873760155654
$ws.Range("T4").Value = 0.1033873760155654

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 79.82725266666667
$ws.Range("H5").Value = 239.481758
$ws.Range("I5").Value = 0.1264482133280045
$ws.Range("J5").Value = 0.1264482133280045
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.446644
$ws.Range("N5").Value = 1.339932
$ws.Range("O5").Value = 0.1823737695100498
$ws.Range("P5").Value = 0.1823737695100498
$ws.Range("Q5").Value = 35.65436344005067
$ws.Range("R5").Value = 320.889270960456
$ws.Range("S5").Value = 0.0230608373124391
$ws.Range("T5").Value = 0.02306083731243909

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 162.5116576666667
$ws.Range("H6").Value = 487.534973
$ws.Range("I6").Value = 0.2574222219914007
$ws.Range("J6").Value = 0.2574222219914007
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.002414333333333
$ws.Range("N6").Value = 6.007243
$ws.Range("O6").Value = 0.8176262304899502
$ws.Range("P6").Value = 0.8176262304899502
$ws.Range("Q6").Value = 325.4156726454932
$ws.Range("R6").Value = 2928.741053809439
$ws.Range("S6").Value = 0.2104751610111761
$ws.Range("T6").Value = 0.2104751610111761

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Agtr2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 162.5116576666667
$ws.Range("H7").Value = 487.534973
$ws.Range("I7").Value = 0.2574222219914007
$ws.Range("J7").Value = 0.2574222219914007
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.446644
$ws.Range("N7").Value = 1.339932
$ws.Range("O7").Value = 0.1823737695100498
$ws.Range("P7").Value = 0.1823737695100498
$ws.Range("Q7").Value = 72.58485682687066
$ws.Range("R7").Value = 653.263711441836
$ws.Range("S7").Value = 0.04694706098022458
$ws.Range("T7").Value = 0.04694706098022458

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Agtr2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 159.7910413333334
$ws.Range("H8").Value = 479.3731240000001
$ws.Range("I8").Value = 0.2531127028358626
$ws.Range("J8").Value = 0.2531127028358626
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.002414333333333
$ws.Range("N8").Value = 6.007243
$ws.Range("O8").Value = 0.8176262304899502
$ws.Range("P8").Value = 0.8176262304899502
$ws.Range("Q8").Value = 319.9678715041259
$ws.Range("R8").Value = 2879.710843537132
$ws.Range("S8").Value = 0.2069515851088093
$ws.Range("T8").Value = 0.2069515851088092

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Agtr2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 159.7910413333334
$ws.Range("H9").Value = 479.3731240000001
$ws.Range("I9").Value = 0.2531127028358626
$ws.Range("J9").Value = 0.2531127028358626
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.446644
$ws.Range("N9").Value = 1.339932
$ws.Range("O9").Value = 0.1823737695100498
$ws.Range("P9").Value = 0.1823737695100498
$ws.Range("Q9").Value = 71.36970986528536
$ws.Range("R9").Value = 642.3273887875681
$ws.Range("S9").Value = 0.04616111772705334
$ws.Range("T9").Value = 0.04616111772705333

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Gnai2"
$ws.Range("C10").Value = "Agtr2"
$ws.Range("D10").Value = "FAPs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 33.73857866666666
$ws.Range("H10").Value = 101.215736
$ws.Range("I10").Value = 0.05344268843173843
$ws.Range("J10").Value = 0.05344268843173842
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.002414333333333
$ws.Range("N10").Value = 6.007243
$ws.Range("O10").Value = 0.8176262304899502
$ws.Range("P10").Value = 0.8176262304899502
$ws.Range("Q10").Value = 67.55861350842754
$ws.Range("R10").Value = 608.0275215758479
$ws.Range("S10").Value = 0.04369614388969116
$ws.Range("T10").Value = 0.04369614388969115

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Gnai2"
$ws.Range("C11").Value = "Agtr2"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 33.73857866666666
$ws.Range("H11").Value = 101.215736
$ws.Range("I11").Value = 0.05344268843173843
$ws.Range("J11").Value = 0.05344268843173842
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.446644
$ws.Range("N11").Value = 1.339932
$ws.Range("O11").Value = 0.1823737695100498
$ws.Range("P11").Value = 0.1823737695100498
$ws.Range("Q11").Value = 15.06913372999467
$ws.Range("R11").Value = 135.622203569952
$ws.Range("S11").Value = 0.00974654454204727
$ws.Range("T11").Value = 0.009746544542047268
